$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B9 held the text "4" — convert it to a genuine number 4
$ws.Range("B9").Value = 4

# Append a new annotation row (row 10)
$ws.Range("A10").Value = "Sunsi Wu"
# Write B10 with a leading apostrophe so Excel keeps it as text "3"
# instead of auto-converting the numeric-looking string to a number,
# then reset the style so no stray quote-prefix formatting is left behind.
$ws.Range("B10").Value = "'3"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "why; not"
$ws.Range("D10").Value = "QSN"
$ws.Range("E10").Value = "RES"
$ws.Range("F10").Value = "a5228610-fe6d-4383-b598-a7c34c3b8714"
$ws.Range("G10").Value = "HyRnez-RW_annotated.xlsx"
$ws.Range("H10").Value = "Why is this result not compared to in Table 1?"
